$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet from "Through 2021-11-13" to "Through 2021-11-14"
$ws.Name = "Through 2021-11-14"

# Update label in A13
$ws.Range("A13").Value = "November (through 11-14)"

# Row 13 updates (November through 11-14)
$ws.Range("C13").Value = 16
$ws.Range("D13").Value = 0.0588
$ws.Range("F13").Value = 32
$ws.Range("G13").Value = 0.0857
$ws.Range("I13").Value = 61
$ws.Range("J13").Value = 0.0161
$ws.Range("L13").Value = 24
$ws.Range("M13").Value = 0.1724
$ws.Range("O13").Value = 19
$ws.Range("P13").Value = 0.1739
$ws.Range("Q13").Value = 3
$ws.Range("R13").Value = 84
$ws.Range("S13").Value = 0.0345
$ws.Range("U13").Value = 96
$ws.Range("V13").Value = 0.0103

# Row 14 updates (Total)
$ws.Range("C14").Value = 242
$ws.Range("D14").Value = 0.12
$ws.Range("F14").Value = 466
$ws.Range("G14").Value = 0.1056
$ws.Range("I14").Value = 710
$ws.Range("J14").Value = 0.0803
$ws.Range("L14").Value = 573
$ws.Range("M14").Value = 0.1102
$ws.Range("O14").Value = 453
$ws.Range("P14").Value = 0.103
$ws.Range("Q14").Value = 57
$ws.Range("R14").Value = 1087
$ws.Range("S14").Value = 0.0498
$ws.Range("U14").Value = 1450
$ws.Range("V14").Value = 0.0578
